$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 2125
$ws.Range("I29").Value = 2000
$ws.Range("J29").Value = 2250
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 6750
$ws.Range("M29").Value = -5719
$ws.Range("N29").Value = -7312
# Row 38
$ws.Range("H38").Value = 421.65216
$ws.Range("I38").Value = 149.44444
$ws.Range("J38").Value = 1401.6
$ws.Range("K38").Value = 448.33332
$ws.Range("L38").Value = 4204.799999999999
$ws.Range("M38").Value = -76.33331999999996
$ws.Range("N38").Value = -4948.799999999999
# Row 58
$ws.Range("H58").Value = 1695.4445
$ws.Range("I58").Value = 433.5
$ws.Range("J58").Value = 2705
$ws.Range("K58").Value = 1300.5
$ws.Range("L58").Value = 8115
$ws.Range("M58").Value = -1150.5
$ws.Range("N58").Value = -8415

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 909.65515
$ws.Range("I2").Value = 702.9524
$ws.Range("J2").Value = 1452.25
$ws.Range("K2").Value = 702.9524
$ws.Range("L2").Value = 1452.25
$ws.Range("M2").Value = -589.9524
$ws.Range("N2").Value = -1678.25
# Row 32
$ws.Range("H32").Value = 8813.857
$ws.Range("I32").Value = 4317.7207
$ws.Range("J32").Value = 23685.691
$ws.Range("K32").Value = 4317.7207
$ws.Range("L32").Value = 23685.691
$ws.Range("M32").Value = -4030.7207
$ws.Range("N32").Value = -24259.691
# Row 61
$ws.Range("H61").Value = 1348.3914
$ws.Range("I61").Value = 1105.2106
$ws.Range("J61").Value = 2503.5
$ws.Range("K61").Value = 1105.2106
$ws.Range("L61").Value = 2503.5
$ws.Range("M61").Value = -893.2106000000001
$ws.Range("N61").Value = -2927.5
# Row 116
$ws.Range("H116").Value = 909.65515
$ws.Range("I116").Value = 702.9524
$ws.Range("J116").Value = 1452.25
$ws.Range("K116").Value = 702.9524
$ws.Range("L116").Value = 1452.25
$ws.Range("M116").Value = 1591.0476
$ws.Range("N116").Value = -6040.25
# Row 132
$ws.Range("H132").Value = 2040.9048
$ws.Range("I132").Value = 1395
$ws.Range("J132").Value = 3203.5334
$ws.Range("K132").Value = 4185
$ws.Range("L132").Value = 9610.600199999999
$ws.Range("M132").Value = -1655
$ws.Range("N132").Value = -14670.6002
# Row 136
$ws.Range("H136").Value = 1348.3914
$ws.Range("I136").Value = 1105.2106
$ws.Range("J136").Value = 2503.5
$ws.Range("K136").Value = 3315.6318
$ws.Range("L136").Value = 7510.5
$ws.Range("M136").Value = -765.6318000000001
$ws.Range("N136").Value = -12610.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 909.65515
$ws.Range("I3").Value = 702.9524
$ws.Range("J3").Value = 1452.25
$ws.Range("K3").Value = 702.9524
$ws.Range("L3").Value = 1452.25
$ws.Range("M3").Value = -588.9524
$ws.Range("N3").Value = -1680.25
# Row 105
$ws.Range("H105").Value = 2008.0952
$ws.Range("I105").Value = 1766.0769
$ws.Range("K105").Value = 1766.0769
$ws.Range("M105").Value = -19.07690000000002
# Row 134
$ws.Range("H134").Value = 2704.3333
$ws.Range("I134").Value = 2304
$ws.Range("J134").Value = 3104.6667
$ws.Range("K134").Value = 6912
$ws.Range("L134").Value = 9314.000100000001
$ws.Range("M134").Value = -4377
$ws.Range("N134").Value = -14384.0001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 59
$ws.Range("H59").Value = 45495
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 45495
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 45495
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -47785
# Row 100
$ws.Range("H100").Value = 20520
$ws.Range("J100").Value = 20520
$ws.Range("L100").Value = 20520
$ws.Range("N100").Value = -22684
# Row 105
$ws.Range("H105").Value = 3564.6667
$ws.Range("I105").Value = 3918.889
$ws.Range("J105").Value = 3033.3333
$ws.Range("K105").Value = 3918.889
$ws.Range("L105").Value = 3033.3333
$ws.Range("M105").Value = -2171.889
$ws.Range("N105").Value = -6527.3333
# Row 132
$ws.Range("H132").Value = 1564.7028
$ws.Range("I132").Value = 1072.4615
$ws.Range("K132").Value = 3217.3845
$ws.Range("M132").Value = -687.3844999999997
# Row 134
$ws.Range("H134").Value = 1896.6666
$ws.Range("I134").Value = 1295.7142
$ws.Range("K134").Value = 3887.1426
$ws.Range("M134").Value = -1352.1426

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 25811
$ws.Range("I57").Value = 17027.5
$ws.Range("J57").Value = 31666.666
$ws.Range("K57").Value = 17027.5
$ws.Range("L57").Value = 31666.666
$ws.Range("M57").Value = -16207.5
$ws.Range("N57").Value = -33306.666
# Row 126
$ws.Range("H126").Value = 1822.2222
$ws.Range("I126").Value = 1600
$ws.Range("J126").Value = 1885.7142
$ws.Range("K126").Value = 4800
$ws.Range("L126").Value = 5657.142599999999
$ws.Range("M126").Value = -2330
$ws.Range("N126").Value = -10597.1426
# Row 132
$ws.Range("H132").Value = 6150.1787
$ws.Range("I132").Value = 8181.8125
$ws.Range("J132").Value = 3441.3333
$ws.Range("K132").Value = 24545.4375
$ws.Range("L132").Value = 10323.9999
$ws.Range("M132").Value = -22015.4375
$ws.Range("N132").Value = -15383.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2095.4
$ws.Range("I40").Value = 2095.4
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2095.4
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1959.4
$ws.Range("N40").ClearContents()
# Row 93
$ws.Range("H93").Value = 6889.7617
$ws.Range("I93").Value = 12400.777
$ws.Range("J93").Value = 2756.5
$ws.Range("K93").Value = 12400.777
$ws.Range("L93").Value = 2756.5
$ws.Range("M93").Value = -11152.777
$ws.Range("N93").Value = -5252.5
# Row 122
$ws.Range("H122").Value = 6394.1816
$ws.Range("I122").Value = 7240.32
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 21720.96
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -19270.96
$ws.Range("N122").Value = -16150

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
# Row 100
$ws.Range("H100").Value = 18183498
$ws.Range("I100").Value = 45454744
$ws.Range("J100").Value = 2666.6667
$ws.Range("K100").Value = 90909488
$ws.Range("L100").Value = 5333.3334
$ws.Range("M100").Value = -90908947
$ws.Range("N100").Value = -6415.3334
# Row 140
$ws.Range("H140").Value = 55738.3
$ws.Range("J140").Value = 55738.3
$ws.Range("L140").Value = 55738.3
$ws.Range("N140").Value = -66098.3
